$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 60, pushing the existing rows 60-142 down to 61-143.
$ws.Rows.Item(60).Insert()

# Populate the newly inserted row 60 with the new record's data.
$ws.Cells.Item(60, 1).Value = 9
$ws.Cells.Item(60, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(60, 3).Value = "Metropolitana"
$ws.Cells.Item(60, 4).Value = 44495
$ws.Cells.Item(60, 5).Value = 13
$ws.Cells.Item(60, 6).Value = 100112017
$ws.Cells.Item(60, 7).Value = "Apio"
$ws.Cells.Item(60, 8).Value = "Americana (o)"
$ws.Cells.Item(60, 9).Value = "Primera"
$ws.Cells.Item(60, 10).Value = 70
$ws.Cells.Item(60, 11).Value = 8000
$ws.Cells.Item(60, 12).Value = 8000
$ws.Cells.Item(60, 13).Value = 8000
$ws.Cells.Item(60, 14).Value = "`$/docena de matas"
$ws.Cells.Item(60, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(60, 16).Value = 1333
$ws.Cells.Item(60, 17).Value = 6
$ws.Cells.Item(60, 18).Value = "Hortaliza"
